# User input added to change basic maze parameters
#
# - "Allow user input" (row 9 of the Networks table, column K) is now DONE,
#   so its status cell M9 flips from TODO (red) to DONE (green).
# - A new Networks TODO is logged on row 10:
#     K10 = "Solve memory issue"
#     L10 = " when adding half a byte"
#     M10 = "TODO" (red)
# - The saved selection moves to C39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "Allow user input" task is complete -> mark DONE (green fill)
$ws.Range("M9").Value = "DONE"
$ws.Range("M9").Interior.Color = 5287936

# Row 10: new Networks task, note, and TODO status (red fill)
$ws.Range("K10").Value = "Solve memory issue"
$ws.Range("L10").Value = " when adding half a byte"
$ws.Range("M10").Value = "TODO"
$ws.Range("M10").Interior.Color = 255

# Restore the workbook's saved cursor position
$ws.Range("C39").Select()
